# Adds two new rows (21 and 22) to the cassette table, reusing the
# border/fill formatting already used by the rows directly above them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Stamp rows 21:22 with the same look (thin box border) as rows 19:20,
# then wipe the copied values/text back out so we can fill in the real
# data for the new rows without the paste clobbering it.
$ws.Range("A19:C20").Copy($ws.Range("A21"))
$ws.Range("A21:C22").ClearContents() | Out-Null

# Row 22 ("BlackBox - 10269277") is typed first so it is registered in
# the shared-strings table ahead of row 21's text, matching the source.
$ws.Range("A22").Value = "BlackBox - 10269277"
$ws.Range("B22").Value = 3
$ws.Range("C22").Value = 72

$muODC = [string]([char]0x252C) + [char]0x00C1 + "ODC"
$ws.Range("A21").Value = $muODC
$ws.Range("B21").Value = 6
$ws.Range("C21").Value = 12

# Leave the selection on the first of the two newly added rows.
$ws.Range("A21").Select() | Out-Null
